# Bharti Airtel Limited - Trading History: add a new, more recent trade.
#
# A new buy (NSE, 2 shares @ 2053, cost 4126.54, CN#252611665409, serial
# date 46062) happened after the most recent existing entry, so it becomes
# the new top data row (row 5). Every existing trade row shifts down by
# one row (old row 5 -> 6, ... old row 19 -> 20); the sheet's used range
# grows from A1:AB19 to A1:AB20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Push the existing data (rows 5-19) down by one row, carrying formatting,
# formulas, etc. along with it.
$ws.Rows.Item(5).Insert()

# The freshly inserted row inherited the bold/shaded header formatting
# from row 4 (Excel's default "format from row above" behavior) - clear
# that so the new row matches the look of the other plain data rows.
$ws.Rows.Item(5).Clear()

# Fill in the new trade's details.
$ws.Range("A5").Value = 46062
$ws.Range("A5").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B5").Value = "NSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2053
$ws.Range("F5").Value = 4126.54
$ws.Range("G5").Value = "CN#252611665409"
$ws.Range("I5").Value = 20.54
$ws.Range("J5").Formula = "=Index!`$C`$2"
